$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) cells we are about to rewrite to Text format
# first. Their values (e.g. "28.331.05", "3.660", "0.02360") must stay as
# literal text -- left as General they would be auto-coerced into numbers
# by Excel, which drops trailing zeros / reinterprets multi-dot strings.
foreach ($r in @(2,3,5,7,8,9,11,12,13,14,15,16,17,19,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,43,44,45,46,47,48,49,50,51)) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = '28.331.05'
$ws.Range("E2").Value = '  -0.54%  '
$ws.Range("D3").Value = '1.810.63'
$ws.Range("E3").Value = '  -0.83%  '
$ws.Range("E4").Value = '  -0.26%  '
$ws.Range("D5").Value = '312.72'
$ws.Range("E5").Value = '  -1.14%  '
$ws.Range("E6").Value = '  -0.21%  '
$ws.Range("D7").Value = '0.5156'
$ws.Range("E7").Value = '  -0.34%  '
$ws.Range("D8").Value = '0.3974'
$ws.Range("E8").Value = '  +2.97%  '
$ws.Range("D9").Value = '0.07831'
$ws.Range("E9").Value = '  -5.59%  '
$ws.Range("E10").Value = '  -1.06%  '
$ws.Range("D11").Value = '41.03'
$ws.Range("E11").Value = '  -2.12%  '
$ws.Range("D12").Value = '6.337'
$ws.Range("E12").Value = '  -0.60%  '
$ws.Range("D13").Value = '1.001'
$ws.Range("E13").Value = '  -0.27%  '
$ws.Range("D14").Value = '20.44'
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = '7.315'
$ws.Range("E15").Value = '  -2.14%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '1.806.27'
$ws.Range("E16").Value = '  -0.94%  '
$ws.Range("D17").Value = '92.47'
$ws.Range("E17").Value = '  -1.77%  '
$ws.Range("E18").Value = '  -3.46%  '
$ws.Range("D19").Value = '0.06568'
$ws.Range("E19").Value = '  -1.03%  '
$ws.Range("D21").Value = '17.31'
$ws.Range("E21").Value = '  -2.80%  '
$ws.Range("D22").Value = '6.009'
$ws.Range("E22").Value = '  -0.65%  '
$ws.Range("D23").Value = '28.349.82'
$ws.Range("E23").Value = '  -0.56%  '
$ws.Range("D24").Value = '11.12'
$ws.Range("E24").Value = '  -3.45%  '
$ws.Range("D25").Value = '2.228'
$ws.Range("E25").Value = '  -0.91%  '
$ws.Range("D26").Value = '160.82'
$ws.Range("E26").Value = '  +0.77%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '20.52'
$ws.Range("E27").Value = '  -2.73%  '
$ws.Range("B28").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C28").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D28").Value = '2.020.02'
$ws.Range("E28").Value = '  -0.65%  '
$ws.Range("D29").Value = '2.422'
$ws.Range("E29").Value = '  +0.82%  '
$ws.Range("D30").Value = '127.74'
$ws.Range("E30").Value = '  +1.38%  '
$ws.Range("D31").Value = '0.1099'
$ws.Range("E31").Value = '  -0.79%  '
$ws.Range("D32").Value = '1.063'
$ws.Range("E32").Value = '  -2.77%  '
$ws.Range("D33").Value = '3.660'
$ws.Range("E33").Value = '  -0.76%  '
$ws.Range("D34").Value = '5.577'
$ws.Range("E34").Value = '  -2.75%  '
$ws.Range("D35").Value = '0.07182'
$ws.Range("E35").Value = '  -4.58%  '
$ws.Range("D36").Value = '9.141'
$ws.Range("E36").Value = '  +4.34%  '
$ws.Range("D37").Value = '0.02360'
$ws.Range("E37").Value = '  -0.19%  '
$ws.Range("D38").Value = '0.2191'
$ws.Range("E38").Value = '  -1.55%  '
$ws.Range("D39").Value = '5.053'
$ws.Range("E39").Value = '  -3.78%  '
$ws.Range("D40").Value = '11.59'
$ws.Range("E40").Value = '  -5.65%  '
$ws.Range("D41").Value = '0.6189'
$ws.Range("E41").Value = '  -3.31%  '
$ws.Range("E42").Value = '  -0.21%  '
$ws.Range("D43").Value = '1.157'
$ws.Range("E43").Value = '  -2.34%  '
$ws.Range("D44").Value = '13.27'
$ws.Range("E44").Value = '  -2.29%  '
$ws.Range("D45").Value = '0.5996'
$ws.Range("E45").Value = '  -3.55%  '
$ws.Range("D46").Value = '1.302'
$ws.Range("E46").Value = '  -6.57%  '
$ws.Range("D47").Value = '3.740'
$ws.Range("E47").Value = '  -1.69%  '
$ws.Range("D48").Value = '125.43'
$ws.Range("E48").Value = '  -1.81%  '
$ws.Range("D49").Value = '1.217'
$ws.Range("E49").Value = '  +0.99%  '
$ws.Range("D50").Value = '1.925'
$ws.Range("E50").Value = '  -4.35%  '
$ws.Range("D51").Value = '0.06835'
$ws.Range("E51").Value = '  -1.80%  '
